$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure numeric-looking price strings in column D stay as text (matches original inlineStr formatting)
$textCells = @("D2","D3","D5","D6","D7","D9","D10","D11","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D26","D27","D28","D29","D30","D33","D34","D35","D36","D38","D40","D42","D46","D47","D48","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '36.441.55'
$ws.Range("E2").Value = '  -2.72%  '
$ws.Range("D3").Value = '1.980.95'
$ws.Range("E3").Value = '  -3.50%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '244.01'
$ws.Range("E5").Value = '  -3.35%  '
$ws.Range("D6").Value = '0.627'
$ws.Range("E6").Value = '  -3.70%  '
$ws.Range("D7").Value = '58.75'
$ws.Range("E7").Value = '  -11.50%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '0.377'
$ws.Range("E9").Value = '  -1.29%  '
$ws.Range("D10").Value = '57.48'
$ws.Range("E10").Value = '  -4.16%  '
$ws.Range("D11").Value = '0.0821'
$ws.Range("E11").Value = '  +7.16%  '
$ws.Range("E12").Value = '  -0.94%  '
$ws.Range("D13").Value = '23.90'
$ws.Range("E13").Value = '  +10.19%  '
$ws.Range("D14").Value = '0.862'
$ws.Range("E14").Value = '  -5.21%  '
$ws.Range("D15").Value = '14.00'
$ws.Range("E15").Value = '  -6.38%  '
$ws.Range("D16").Value = '2.271.21'
$ws.Range("E16").Value = '  -3.56%  '
$ws.Range("D17").Value = '5.46'
$ws.Range("E17").Value = '  -2.26%  '
$ws.Range("D18").Value = '1.983.13'
$ws.Range("E18").Value = '  -3.37%  '
$ws.Range("D19").Value = '36.323.67'
$ws.Range("E19").Value = '  -2.60%  '
$ws.Range("D20").Value = '70.58'
$ws.Range("E20").Value = '  -4.19%  '
$ws.Range("D21").Value = '0.0₃0861'
$ws.Range("E21").Value = '  -1.84%  '
$ws.Range("D22").Value = '5.32'
$ws.Range("E22").Value = '  -2.28%  '
$ws.Range("D23").Value = '234.42'
$ws.Range("E23").Value = '  -2.39%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("E25").Value = '  -2.70%  '
$ws.Range("D26").Value = '2.31'
$ws.Range("E26").Value = '  -3.82%  '
$ws.Range("D27").Value = '10.13'
$ws.Range("E27").Value = '  +3.25%  '
$ws.Range("D28").Value = '161.68'
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").Value = '19.85'
$ws.Range("E29").Value = '  -1.07%  '
$ws.Range("D30").Value = '0.127'
$ws.Range("E30").Value = '  +7.98%  '
$ws.Range("E31").Value = '  -1.60%  '
$ws.Range("E32").Value = '  -0.23%  '
$ws.Range("D33").Value = '4.90'
$ws.Range("E33").Value = '  -6.85%  '
$ws.Range("D34").Value = '0.0633'
$ws.Range("E34").Value = '  +1.50%  '
$ws.Range("D35").Value = '4.42'
$ws.Range("E35").Value = '  -6.00%  '
$ws.Range("D36").Value = '6.28'
$ws.Range("E36").Value = '  +3.88%  '
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").Value = '2.27'
$ws.Range("E38").Value = '  -7.23%  '
$ws.Range("E39").Value = '  -4.44%  '
$ws.Range("D40").Value = '3.09'
$ws.Range("E40").Value = '  +2.70%  '
$ws.Range("E41").Value = '  +0.42%  '
$ws.Range("D42").Value = '0.0963'
$ws.Range("E42").Value = '  -6.72%  '
$ws.Range("E43").Value = '  -3.49%  '
$ws.Range("E44").Value = '  -2.74%  '
$ws.Range("E45").Value = '  -4.72%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '92.64'
$ws.Range("E46").Value = '  -3.30%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '16.22'
$ws.Range("E47").Value = '  -4.78%  '
$ws.Range("D48").Value = '7.56'
$ws.Range("E48").Value = '  -5.20%  '
$ws.Range("D49").Value = '1.375.94'
$ws.Range("E49").Value = '  -3.28%  '
$ws.Range("E50").Value = '  -3.29%  '
$ws.Range("D51").Value = '45.05'
$ws.Range("E51").Value = '  -3.63%  '
